$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New draw result row appended by the daily auto-update job.
# Leading apostrophes force Date-like / numeric-looking strings to be
# stored as text (matching the existing column formatting), then
# ClearFormats() strips the transient "quote prefix" cell style Excel
# applies for that so the cell keeps the workbook's default style.

$ws.Range("A11").Value = "'2025-09-27"
$ws.Range("A11").ClearFormats()

$ws.Range("B11").Value = "Pick 4"

$ws.Range("C11").Value = "'250927"
$ws.Range("C11").ClearFormats()

$ws.Range("D11").Value = "5-6-0-1"

$ws.Range("E11").Value = "2025-09-27T21:34:59.820+04:00"
